# Add a second (data) row below the existing header row on the
# Master_Linking_Log sheet.
#
# The new values are record-like strings ("123", "1234567",
# "2017-11-05", "18:48:25.326411", "18", ...) that Excel would normally
# auto-detect as numbers/dates. We want them stored as plain text (shared
# strings), matching the source data. To force text entry without Excel's
# automatic number/date conversion, the target range is temporarily
# switched to a text number format, the values are written, and then the
# temporary formatting is cleared again (ClearFormats) so the cells keep
# their default style while their stored type remains text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRow = $ws.Range("A2:H2")
$dataRow.NumberFormat = "@"

$ws.Range("A2").Value = "123"
$ws.Range("B2").Value = "1234567"
$ws.Range("C2").Value = "2017-11-05"
$ws.Range("D2").Value = "18:48:25.326411"
$ws.Range("E2").Value = "abc"
$ws.Range("F2").Value = "18"
$ws.Range("G2").Value = "m"
$ws.Range("H2").Value = "asb"

$dataRow.ClearFormats()
